$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the "lag" time window for temperature/precipitation covariates in
# favor of the "ann" (annual) window.
$ws.Range("B14").Value = "temp_ann_cont_scale_clst"
$ws.Range("B15").Value = "precip_ann_cont_scale_clst"

# Reflect the last-edited cell in the sheet view selection.
$ws.Range("B15").Select()
